$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Valid_Login")

# Drop every existing hyperlink on the sheet so stale relationship
# ids/targets don't survive the row reshuffle below.
$ws.Hyperlinks.Delete()

# The first four data rows (opensourcecms x2, demo1, brnwinjit) are gone;
# deleting them shifts the surekha / anand / Nasreen rows up to rows 2-4.
$ws.Rows("2:5").Delete()

# Row 2 (was surekha, row 6) is now the "Y" / passing record, with a new
# password.
$ws.Range("A2").Value = "Y"
$ws.Range("C2").Value = "Surekha@123123"

# Row 4 (was Nasreen, row 8) is no longer the "Y" record.
$ws.Range("A4").Value = "N"

# Re-create the hyperlinks against the new row layout.
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:surekha.jadhav@ram.co.za")
$ws.Range("B2").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:anand.maurya@ram.co.za")
$ws.Range("B3").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:Surekha@123123")
$ws.Range("C2").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:Nasreen.khan@ram.co.za")
$ws.Range("B4").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:Nasreen@123")
$ws.Range("C4").Style = "Hyperlink"

# Match the author's final selection on this sheet.
$ws.Range("B12").Select()
